$wb = $excel.ActiveWorkbook
$taxSheet = $wb.Worksheets.Item("TAXONOMY")

# Create the new worksheet positioned after TAXONOMY
$newSheet = $wb.Worksheets.Add([System.Reflection.Missing]::Value, $taxSheet)
$newSheet.Name = "EXTENDED_2025_JULY"

# Copy header row (A1:L1) from TAXONOMY to preserve header styling
$taxSheet.Range("A1:L1").Copy($newSheet.Range("A1:L1"))

# Row 2: copy formatting template from TAXONOMY row 5
$taxSheet.Range("A5:L5").Copy($newSheet.Range("A2:L2"))
# Row 3: copy formatting template from TAXONOMY row 5
$taxSheet.Range("A5:L5").Copy($newSheet.Range("A3:L3"))
# Row 4: copy formatting template from TAXONOMY row 2
$taxSheet.Range("A2:L2").Copy($newSheet.Range("A4:L4"))
# Row 5: copy formatting template from TAXONOMY row 6
$taxSheet.Range("A6:L6").Copy($newSheet.Range("A5:L5"))
# Row 6: copy formatting template from TAXONOMY row 5
$taxSheet.Range("A5:L5").Copy($newSheet.Range("A6:L6"))
# Row 7: copy formatting template from TAXONOMY row 2
$taxSheet.Range("A2:L2").Copy($newSheet.Range("A7:L7"))
# Row 8: copy formatting template from TAXONOMY row 2
$taxSheet.Range("A2:L2").Copy($newSheet.Range("A8:L8"))
# Row 9: copy formatting template from TAXONOMY row 6
$taxSheet.Range("A6:L6").Copy($newSheet.Range("A9:L9"))

# --- Row 2 data ---
$newSheet.Range("A2").Value = 'icozAutomatedCodeReview2025'
$v1 = @'
```bibtex
@preprint{icozAutomatedCodeReview2025,
      title={{Automated Code Review Using Large Language Models with Symbolic Reasoning}}, 
      author={Busra Icoz and Goksel Biricik},
      year={2025},
      eprint={2507.18476},
      archivePrefix={arXiv}}```
'@
$newSheet.Range("B2").Value = $v1
$newSheet.Range("C2").Value = 'https://arxiv.org/abs/2507.18476'
$newSheet.Range("G2").Value = 'Binary'
$newSheet.Range("H2").Value = 'Reasoning'
$newSheet.Range("I2").Value = 'Code LMs'
$newSheet.Range("J2").Value = 'Full-Parameter Fine-Tuning, Few-Shot'
$newSheet.Range("K2").Value = 'Raw'
$newSheet.Range("L2").Value = 'Devign'
$newSheet.Range("E2").Value = 45862
$newSheet.Range("F2").Value = 45862
$newSheet.Range("D2").ClearContents()

# --- Row 3 data ---
$newSheet.Range("A3").Value = 'zhangBountyBenchDollarImpact2025'
$v2 = @'
```bibtex
@preprint{zhangBountyBenchDollarImpact2025,
      title={{BountyBench: Dollar Impact of AI Agent Attackers and Defenders on Real-World Cybersecurity Systems}}, 
      author={Andy K. Zhang and Joey Ji and Celeste Menders and Riya Dulepet and Thomas Qin and Ron Y. Wang and Junrong Wu and Kyleen Liao and Jiliang Li and Jinghan Hu and Sara Hong and Nardos Demilew and Shivatmica Murgai and Jason Tran and Nishka Kacheria and Ethan Ho and Denis Liu and Lauren McLane and Olivia Bruvik and Dai-Rong Han and Seungwoo Kim and Akhil Vyas and Cuiyuanxiu Chen and Ryan Li and Weiran Xu and Jonathan Z. Ye and Prerit Choudhary and Siddharth M. Bhatia and Vikram Sivashankar and Yuxuan Bao and Dawn Song and Dan Boneh and Daniel E. Ho and Percy Liang},
      year={2025},
      eprint={2505.15216},
      archivePrefix={arXiv}}```
'@
$newSheet.Range("B3").Value = $v2
$newSheet.Range("C3").Value = 'https://arxiv.org/abs/2505.15216'
$newSheet.Range("G3").Value = 'Multi-Class, Vulnerability-Specific'
$newSheet.Range("H3").Value = 'Repair, Exploit'
$newSheet.Range("I3").Value = 'General LMs, Code LMs'
$newSheet.Range("J3").Value = 'Agentic, Cot'
$newSheet.Range("K3").Value = 'Prompt'
$newSheet.Range("L3").Value = 'Custom'
$newSheet.Range("E3").Value = 45798
$newSheet.Range("F3").Value = 45848
$newSheet.Range("D3").ClearContents()

# --- Row 4 data ---
$newSheet.Range("A4").Value = 'liCLeVeRMultimodalContrastive'
$v3 = @'
```bibtex
@inproceedings{liCLeVeRMultimodalContrastive,
    title = {{CLeVeR: Multi-modal Contrastive Learning for Vulnerability Code Representation}},
    author = {Li, Jiayuan  and Cui, Lei  and Zhao, Sen  and Yang, Yun  and Li, Lun  and Zhu, Hongsong},
    booktitle = {Findings of the Association for Computational Linguistics (ACL)},
    year = {2025},
    address = {Vienna, Austria},
    publisher = {ACL},
    doi = {10.18653/v1/2025.findings-acl.414},
    pages = {7940--7951}
}```
'@
$newSheet.Range("B4").Value = $v3
$newSheet.Range("C4").Value = 'https://aclanthology.org/2025.findings-acl.414/'
$newSheet.Range("D4").Value = 'https://github.com/yoimiya-nlp/CLeVeR'
$newSheet.Range("G4").Value = 'Binary, Multi-Class'
$newSheet.Range("I4").Value = 'General LMs, Code LMs'
$newSheet.Range("J4").Value = 'Contrastive Learning, Pre-Training, PEFT (linear probing), Adapter-Tuning'
$newSheet.Range("K4").Value = 'Raw, Structure-Aware'
$newSheet.Range("L4").Value = 'SARD, SynData, Devign, Reveal, Custom, VCLData'
$newSheet.Range("E4").Value = 45839
$newSheet.Range("F4").Value = 45840
$newSheet.Range("H4").ClearContents()

# --- Row 5 data ---
$newSheet.Range("A5").Value = 'sunHgtJITJustinTimeVulnerability2025'
$v4 = @'
```bibtex
@article{sunHgtJITJustinTimeVulnerability2025,
  author={Sun, Xiaobing and Zhou, Mingxuan and Cao, Sicong and Wu, Xiaoxue and Bo, Lili and Wu, Di and Li, Bin and Xiang, Yang},
  journal={IEEE Transactions on Dependable and Secure Computing (TDSC)}, 
  title={{HgtJIT: Just-in-Time Vulnerability Detection Based on Heterogeneous Graph Transformer}}, 
  year={2025},  pages={1-17},
  doi={10.1109/TDSC.2025.3586669}}```
'@
$newSheet.Range("B5").Value = $v4
$newSheet.Range("C5").Value = 'https://ieeexplore.ieee.org/abstract/document/11072308'
$newSheet.Range("G5").Value = 'Binary'
$newSheet.Range("I5").Value = 'Hybrid - GNN'
$newSheet.Range("J5").Value = 'None'
$newSheet.Range("K5").Value = 'Structure-Aware'
$newSheet.Range("L5").Value = 'CodeJIT'
$newSheet.Range("E5").Value = 45845
$newSheet.Range("F5").Value = 45846
$newSheet.Range("D5").ClearContents()
$newSheet.Range("H5").ClearContents()

# --- Row 6 data ---
$newSheet.Range("A6").Value = 'simoniImprovingLLMReasoning2025'
$v5 = @'
```bibtex
@preprint{simoniImprovingLLMReasoning2025,
      title={{Improving LLM Reasoning for Vulnerability Detection via Group Relative Policy Optimization}}, 
      author={Marco Simoni and Aleksandar Fontana and Giulio Rossolini and Andrea Saracino},
      year={2025},
      eprint={2507.03051},
      archivePrefix={arXiv}
}```
'@
$newSheet.Range("B6").Value = $v5
$newSheet.Range("C6").Value = 'https://arxiv.org/abs/2507.03051'
$newSheet.Range("G6").Value = 'Binary'
$newSheet.Range("H6").Value = 'Reasoning'
$newSheet.Range("I6").Value = 'General LMs'
$newSheet.Range("J6").Value = 'Zero-Shot, CoT, Full-Parameter Fine-Tuning'
$newSheet.Range("K6").Value = 'Prompt, Raw'
$newSheet.Range("L6").Value = 'Big-Vul, DiverseVul, CleanVul'
$newSheet.Range("E6").Value = 45841
$newSheet.Range("F6").Value = 45842
$newSheet.Range("D6").ClearContents()

# --- Row 7 data ---
$newSheet.Range("A7").Value = 'lekssaysLLMxCPGContextAwareVulnerability2025'
$v6 = @'
```bibtex
@preprint{lekssaysLLMxCPGContextAwareVulnerability2025,
      title={{LLMxCPG: Context-Aware Vulnerability Detection Through Code Property Graph-Guided Large Language Models}}, 
      author={Ahmed Lekssays and Hamza Mouhcine and Khang Tran and Ting Yu and Issa Khalil},
      year={2025},
      eprint={2507.16585},
      archivePrefix={arXiv}}```
'@
$newSheet.Range("B7").Value = $v6
$newSheet.Range("C7").Value = 'https://arxiv.org/abs/2507.16585'
$newSheet.Range("D7").Value = 'https://github.com/qcri/llmxcpg; https://zenodo.org/records/15614095'
$newSheet.Range("G7").Value = 'Binary'
$newSheet.Range("I7").Value = 'General LMs'
$newSheet.Range("J7").Value = 'Low-Rank Decomposition'
$newSheet.Range("K7").Value = 'Structure-Aware, Prompt'
$newSheet.Range("L7").Value = 'FormAI, PrimeVul, SVEN, ReposVul'
$newSheet.Range("E7").Value = 45860
$newSheet.Range("F7").Value = 45860
$newSheet.Range("H7").ClearContents()

# --- Row 8 data ---
$newSheet.Range("A8").Value = 'liOutDistributionOut2025'
$v7 = @'
```bibtex
@preprint{liOutDistributionOut2025,
      title={{Out of Distribution, Out of Luck: How Well Can LLMs Trained on Vulnerability Datasets Detect Top 25 CWE Weaknesses?}}, 
      author={Yikun Li and Ngoc Tan Bui and Ting Zhang and Martin Weyssow and Chengran Yang and Xin Zhou and Jinfeng Jiang and Junkai Chen and Huihui Huang and Huu Hung Nguyen and Chiok Yew Ho and Jie Tan and Ruiyin Li and Yide Yin and Han Wei Ang and Frank Liauw and Eng Lieh Ouh and Lwin Khin Shar and David Lo},
      year={2025},
      eprint={2507.21817},
      archivePrefix={arXiv}
```
'@
$newSheet.Range("B8").Value = $v7
$newSheet.Range("C8").Value = 'https://arxiv.org/abs/2507.21817'
$newSheet.Range("D8").Value = 'https://github.com/yikun-li/TitanVul-BenchVul'
$newSheet.Range("G8").Value = 'Multi-Class'
$newSheet.Range("I8").Value = 'General LMs, Code LMs'
$newSheet.Range("J8").Value = 'Full-Parameter Fine-Tuning'
$newSheet.Range("K8").Value = 'Raw'
$newSheet.Range("L8").Value = 'Custom, BenchVul, TitanVul'
$newSheet.Range("E8").Value = 45867
$newSheet.Range("F8").Value = 45883
$newSheet.Range("H8").ClearContents()

# --- Row 9 data ---
$newSheet.Range("A9").Value = 'liRevisitingPretrainedLanguage2025'
$v8 = @'
```bibtex
@preprint{liRevisitingPretrainedLanguage2025,
      title={{Revisiting Pre-trained Language Models for Vulnerability Detection}}, 
      author={Youpeng Li and Weiliang Qi and Xuyu Wang and Fuxun Yu and Xinda Wang},
      year={2025},
      eprint={2507.16887},
      archivePrefix={arXiv}
```
'@
$newSheet.Range("B9").Value = $v8
$newSheet.Range("C9").Value = 'https://arxiv.org/abs/2507.16887'
$newSheet.Range("G9").Value = 'Binary, Multi-Class'
$newSheet.Range("I9").Value = 'General LMs, Code LMs'
$newSheet.Range("J9").Value = 'Zero-Shot, Few-Shot, Full-Parameter Fine-Tuning, Low-Rank Decomposition'
$newSheet.Range("K9").Value = 'Structure-Aware, Prompt'
$newSheet.Range("L9").Value = 'PrimeVul, Custom'
$newSheet.Range("E9").Value = 45860
$newSheet.Range("F9").Value = 45860
$newSheet.Range("D9").ClearContents()
$newSheet.Range("H9").ClearContents()

# Create the table (ListObject) over the new data, matching Tabelle132 naming
$lo = $newSheet.ListObjects.Add(1, $newSheet.Range("A1:L9"), $null, 1)
$lo.Name = "Tabelle132"
$lo.TableStyle = "TableStyleLight9"

# Set the active selection on the new sheet
$newSheet.Range("B4").Select()

# Restore TAXONOMY as the active sheet and update its selection
$taxSheet.Activate()
$taxSheet.Range("A24").Select()
